# Corrected file endings: orig_filename column (A) listed the source
# annotation files with the old ".xmi" extension; the files were renamed
# to ".tsv", so update the cell values to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$renames = @{
    "A2"  = "rwz_full_digbib_1136.tsv"
    "A3"  = "rwz_full_digbib_1199.tsv"
    "A4"  = "rwz_full_digbib_1208.tsv"
    "A5"  = "rwz_full_digbib_1340.tsv"
    "A6"  = "rwz_full_digbib_2473.tsv"
    "A7"  = "rwz_full_digbib_2632.tsv"
    "A8"  = "rwz_full_digbib_2649.tsv"
    "A9"  = "rwz_full_digbib_3152.tsv"
    "A10" = "rwz_full_digbib_3153.tsv"
    "A11" = "rwz_full_digbib_3228.tsv"
    "A12" = "rwz_full_digbib_5000.tsv"
    "A13" = "rwz_full_digbib_5001.tsv"
    "A14" = "rwz_full_digbib_5002.tsv"
    "A15" = "rwz_full_digbib_5003.tsv"
    "A16" = "rwz_full_digbib_5004.tsv"
    "A17" = "rwz_full_digbib_5005.tsv"
    "A18" = "rwz_full_digbib_5006.tsv"
    "A19" = "rwz_full_digbib_5007.tsv"
}

foreach ($cellRef in $renames.Keys) {
    $ws.Range($cellRef).Value = $renames[$cellRef]
}
